$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$oldLine1 = "✅ 1000 Bs = 14.95 = 62152.47 pesos"
$oldLine2 = "✅ 62152.47 pesos = 14.94 = 980.77 Bs"
$newLine1 = "✅ 1000 Bs = 14.58 = 59932.36 pesos"
$newLine2 = "✅ 59932.36 pesos = 14.47 = 969.2 Bs"

$text = $ws.Range("A1").Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws.Range("A1").Value = $text

$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 68.59399999999999
$wsTasas.Range("O10").Value = 4111
$wsTasas.Range("N12").Value = 4143
$wsTasas.Range("O12").Value = 66.999
